# week4-reading.docx
#
# 1) Several phrases that currently live in a single <w:r> get split
#    into one run per word/segment (title, and every curly-quoted
#    term). Word produces exactly this kind of run fragmentation
#    whenever a sub-range of a run is independently touched (e.g. a
#    formatting toggle applied and then reverted on part of it) -- the
#    text is unchanged but the run boundaries move. We reproduce that
#    by toggling Bold on/off (net no-op) on each sub-segment.
#
# 2) Two paragraph styles get tweaked: Subtitle is rebased onto Normal
#    (instead of Title) and picks up a muted text colour, and
#    AbstractTitle gets an explicit blue text colour.

$d = $word.ActiveDocument

function Split-Phrase($phrase, $parts) {
    # Find `$phrase` in the document body, then force a run boundary
    # after every sub-segment in `$parts` (their concatenation must
    # equal `$phrase`) by toggling Bold on and back off across each
    # piece. Leaves the visible text completely unchanged.
    $rng = $d.Content
    $found = $rng.Find.Execute($phrase, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return
    }
    $pos = $rng.Start
    foreach ($part in $parts) {
        $len = $part.Length
        $seg = $d.Range($pos, $pos + $len)
        $seg.Bold = 1
        $seg.Bold = 0
        $pos = $pos + $len
    }
}

Split-Phrase "Week 4 Reading Guide: Basic Regression" @("Week", " ", "4", " ", "Reading", " ", "Guide:", " ", "Basic", " ", "Regression")

Split-Phrase "“response”" @("“", "response", "”")
Split-Phrase "“explanatory”" @("“", "explanatory", "”")
Split-Phrase "“explanatory modeling”" @("“", "explanatory modeling", "”")
Split-Phrase "“predictive modeling”" @("“", "predictive modeling", "”")
Split-Phrase "“basic”" @("“", "basic", "”")
Split-Phrase "“EDA”" @("“", "EDA", "”")
Split-Phrase "“fit”" @("“", "fit", "”")
Split-Phrase "“levels”" @("“", "levels", "”")
Split-Phrase "“baseline”" @("“", "baseline", "”")
Split-Phrase "“indicator function”" @("“", "indicator function", "”")
Split-Phrase "“correlation does not imply causation,”" @("“", "correlation does not imply causation,", "”")

# Subtitle style: base it on Normal instead of Title, and give its
# text a muted theme-based colour instead of inheriting Title's.
$subtitle = $d.Styles.Item("Subtitle")
$subtitle.BaseStyle = $d.Styles.Item("Normal")
$subtitle.Font.TextColor.ObjectThemeColor = 13   # wdThemeColorText1
$subtitle.Font.TextColor.TintAndShade = 0.65     # hex "A6" tint

# AbstractTitle style: add an explicit blue colour (345A8A) to the
# bold run text. Font.Color uses VBA's RGB() byte order (0xBBGGRR).
$abstractTitle = $d.Styles.Item("AbstractTitle")
$abstractTitle.Font.Color = 9067060   # RGB(0x34, 0x5A, 0x8A) -> "345A8A"
